$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Blåsippa) and row 11 (Grönpyrola) have had their species-specific
# data swapped. Columns C,D,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY are
# identical between the two rows, so only A,B,E,F,G,H,I,J,K,Q,R need to be
# exchanged.
$cols = @("A","B","E","F","G","H","I","J","K","Q","R")

# Column I holds a numeric-looking count ("35") that is stored as TEXT in
# the source data rather than as a number. Reading it back via COM
# (Value2/Text/Formula) always normalises it to a numeric variant, so the
# text-ness can't be detected from the read side - special-case it here by
# column instead, and quote-prefix it on write so it round-trips as text.
$textLikeNumberCols = @("I")

foreach ($col in $cols) {
    $addr10 = $col + "10"
    $addr11 = $col + "11"
    $v10 = $ws.Range($addr10).Value2
    $v11 = $ws.Range($addr11).Value2

    if ($textLikeNumberCols -contains $col) {
        if ("" + $v11 -ne "") {
            $ws.Range($addr10).Value2 = "'" + $v11
        } else {
            $ws.Range($addr10).Value2 = $v11
        }
        if ("" + $v10 -ne "") {
            $ws.Range($addr11).Value2 = "'" + $v10
        } else {
            $ws.Range($addr11).Value2 = $v10
        }
    } else {
        $ws.Range($addr10).Value2 = $v11
        $ws.Range($addr11).Value2 = $v10
    }
}
